$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (2013)
$ws.Range("B2").Value = 1727120.135999999
$ws.Range("E2").Value = 349996.8559999999
$ws.Range("F2").Value = 474310.8660000003
$ws.Range("G2").Value = 1361191.480000003
$ws.Range("H2").Value = 316120.0299999998

# Row 3 (2014)
$ws.Range("B3").Value = 1903613.317000007
$ws.Range("E3").Value = 234735.071
$ws.Range("F3").Value = 648629.6440000001
$ws.Range("G3").Value = 1463516.151999999
$ws.Range("H3").Value = 380410.7770000005

# Row 4 (2015)
$ws.Range("B4").Value = 1977125.756000004
$ws.Range("D4").Value = 92324.52199999998
$ws.Range("F4").Value = 470395.9069999999
$ws.Range("G4").Value = 1329979.525000002
$ws.Range("H4").Value = 501160.7859999991

# Row 5 (2016)
$ws.Range("B5").Value = 2295864.517
$ws.Range("E5").Value = 378440.794
$ws.Range("F5").Value = 459383.149
$ws.Range("G5").Value = 1639697.842000007
$ws.Range("H5").Value = 434037.784

# Row 6 (2017)
$ws.Range("B6").Value = 2451251.181999999
$ws.Range("C6").Value = 152105.5440000001
$ws.Range("E6").Value = 289694.662
$ws.Range("G6").Value = 1783334.905000001
$ws.Range("H6").Value = 790149.3520000008

# Row 7 (2018)
$ws.Range("B7").Value = 2272601.563
$ws.Range("C7").Value = 83341.09199999999
$ws.Range("E7").Value = 418200.9099999995
$ws.Range("G7").Value = 1936112.166000004
$ws.Range("H7").Value = 556489.875999998
